$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 - Enterprises density (per 1000 people): Micro and MSMEs columns
$ws.Range("B13").Value = "'69.58"
$ws.Range("D13").Value = "'80.69"

# Row 14 - Employment (% of total): Micro, SMEs, MSMEs columns
$ws.Range("B14").Value = "'24.18"
$ws.Range("C14").Value = "'44.87"
$ws.Range("D14").Value = "'69.05"

# Row 16 - Enterprises (% of total): Micro, SMEs, MSMEs columns
$ws.Range("B16").Value = "'85.95"
$ws.Range("C16").Value = "'13.71"
$ws.Range("D16").Value = "'99.66"
